$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("13:13").Insert()
$ws.Range("A13").Value = "c_sand"
$ws.Range("B13").Value = "EA 23.112"
"done"
